$d = $word.ActiveDocument

# Read the full document OOXML (this includes <w:sdt> wrappers, unlike
# Range.WordOpenXML scoped to a content control's inner range).
$xml = $d.Content.WordOpenXML

$old = '<w:sdtPr><w:id w:val="-1388647104"/><w:docPartObj><w:docPartGallery w:val="Table of Contents"/><w:docPartUnique/></w:docPartObj></w:sdtPr><w:sdtEndPr><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorBidi"/><w:b/><w:bCs/><w:noProof/><w:color w:val="auto"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:sdtEndPr>'

$new = '<w:sdtPr><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorBidi"/><w:color w:val="auto"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:id w:val="-1388647104"/><w:docPartObj><w:docPartGallery w:val="Table of Contents"/><w:docPartUnique/></w:docPartObj></w:sdtPr><w:sdtEndPr><w:rPr><w:b/><w:bCs/><w:noProof/></w:rPr></w:sdtEndPr>'

if ($xml.IndexOf($old) -lt 0) {
    throw "Expected sdtPr/sdtEndPr block for the Table of Contents content control was not found."
}

$xml = $xml.Replace($old, $new)

$d.Content.WordOpenXML = $xml
